# Add a new service event (row 18) to the Card17 sheet and backfill the
# "nan" placeholders in row 17 that precede it, per the Machine Service
# Lookup workbook's existing convention for blank measurement columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card17")

# --- Row 17: fill the previously-empty measurement columns with "nan" ---
$row17Cols = @(2,3,4,5,6,7,8,9,10,11,13)  # B,C,D,E,F,G,H,I,J,K,M
foreach ($col in $row17Cols) {
    $ws.Cells.Item(17, $col).Value2 = "nan"
}

# --- Row 18: brand-new service record ---
$ws.Cells.Item(18, 1).Value2 = "17"
# B18:K18 stay blank (matches the sheet's "no measurement" convention)
$ws.Cells.Item(18, 12).Value2 = "22\10\2025"
# M18 stays blank
$ws.Cells.Item(18, 14).Value2 = "تم تغيير الجرائد الخلفيه (1_5_8) ومعايرتها"
$ws.Cells.Item(18, 15).Value2 = "م\محمد عبدالله"
